# "Generate Report for Handoff" - refresh the handoff status/report:
#   - Overview sheet: status moves from "Handed back: in sync with en-US"
#     to "Ready for handoff", and the Latest HO Xliff Generate Date is
#     bumped to the new handoff timestamp.
#   - zh-cn sheet: Status mirrors the Overview status change, and Latest
#     Handoff Datetime is bumped to the new handoff timestamp.
#   - de-de sheet: Status mirrors the Overview status change, and Latest
#     Handoff Datetime (same text as the Overview's new generate date) is
#     bumped too.
#   - The status/date columns also got narrower on handoff regeneration.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 14:36:59"

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 14:36:46"

# --- de-de sheet (same generate-date text as Overview!G2) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-07 14:36:59"

# --- Column width tweaks on the refreshed status/date columns ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
